$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-7 (A,B,C columns) are rotated: row4's data moves to row7,
# while rows 5,6,7 shift up to 4,5,6.
$data = @(
    @(4, "Thangkas  - groß", "große Thangka"),
    @(5, "Thangkas  - groß", "große Thangka 2"),
    @(6, "Thangkas  - klein", "kleine thangka"),
    @(3, "Malas  - klein", "kleine mala")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 4 + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Update selection: activeCell A4, selected range A4:XFD7 (entire rows 4:7)
$ws.Range("A4:XFD7").Select()
